$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write the brand-new unit/label strings first, in the order that reproduces
#     the target shared-string table layout (Kilograms, Centimeters, Lamina*,
#     Primary*, Thallus mass) ---
$ws.Range("C2").Value  = "Kilograms (kg)"
$ws.Range("C3").Value  = "Centimeters (cm)"

$ws.Range("B12").Value = "Lamina weight"
$ws.Range("B13").Value = "Lamina length"
$ws.Range("B14").Value = "Lamina thickness"

$ws.Range("B3").Value  = "Primary length"
$ws.Range("B4").Value  = "Primary width"

$ws.Range("B19").Value = "Thallus mass"

# --- Fill in the remaining cells that reuse already-created strings ---
$ws.Range("C4").Value  = "Centimeters (cm)"
$ws.Range("C5").Value  = "Centimeters (cm)"
$ws.Range("C6").Value  = "Kilograms (kg)"
$ws.Range("C7").Value  = "Centimeters (cm)"
$ws.Range("C10").Value = "Centimeters (cm)"
$ws.Range("C11").Value = "Centimeters (cm)"

$ws.Range("C12").Value = "Kilograms (kg)"
$ws.Range("C13").Value = "Centimeters (cm)"
$ws.Range("C15").Value = "Kilograms (kg)"
$ws.Range("C16").Value = "Centimeters (cm)"

# C13, C15, C16 switch to the unwrapped/bordered style used in rows 2-11 (style index 3)
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add two new rows (19, 20) for Laminaria pallida, copying formatting from row 18 ---
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A19").Value = "Laminaria pallida"
$ws.Range("C19").Value = "Grams (g)"

$ws.Range("A20").Value = "Laminaria pallida"
$ws.Range("B20").Value = "Total length"
$ws.Range("C20").Value = "Centimeters (cm)"

# --- Update selection ---
$ws.Range("E17").Select()
